{"js": "// Rewrite the \"Impact\" bullet list under \"KEY ACHIEVEMENTS AND IMPACT\"\n// from 6 job-duty style bullets into 4 impact-focused accomplishment\n// statements, per the commit: \"Fix Key Achievements to use proper\n// accomplishment statements\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" section heading. There is\n// only one such heading in the document, so this anchors us reliably\n// even though some bullet text is duplicated elsewhere (e.g. under\n// \"Partner - Siege Analytics\").\nconst items = paragraphs.items;\nlet sectionHeadingIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionHeadingIndex = i;\n    break;\n  }\n}\nif (sectionHeadingIndex === -1) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\n// Collect the bullet paragraphs that belong to this section: everything\n// after the heading (and its \"Impact\" sub-heading) up to the next\n// heading-styled paragraph (Heading 1/2/3).\nconst isHeadingStyle = (style) => /^Heading/i.test(style || \"\");\n\nlet bulletStart = sectionHeadingIndex + 1;\n// Skip an optional sub-heading paragraph (e.g. \"Impact\").\nif (bulletStart < items.length && isHeadingStyle(items[bulletStart].style)) {\n  bulletStart++;\n}\n\nlet bulletEnd = bulletStart; // exclusive\nwhile (bulletEnd < items.length && !isHeadingStyle(items[bulletEnd].style)) {\n  bulletEnd++;\n}\n\nconst bulletParagraphs = items.slice(bulletStart, bulletEnd);\n\nconst newBullets = [\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\u2022 $4.7M savings enabled nonprofit access\",\n  \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"\u2022 178% accuracy improvement in racial classification algorithms\",\n];\n\nif (bulletParagraphs.length < newBullets.length) {\n  throw new Error(\"Fewer existing bullet paragraphs than replacement bullets\");\n}\n\n// Overwrite the text of the first N bullet paragraphs in place so their\n// formatting/paragraph identity is preserved, then delete any leftover\n// paragraphs beyond the new bullet count.\nfor (let i = 0; i < newBullets.length; i++) {\n  bulletParagraphs[i].insertText(newBullets[i], Word.InsertLocation.replace);\n}\nfor (let i = newBullets.length; i < bulletParagraphs.length; i++) {\n  bulletParagraphs[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"Impact\" bullet list under \"KEY ACHIEVEMENTS AND IMPACT\"\n# from 6 job-duty style bullets into 4 impact-focused accomplishment\n# statements, per the commit: \"Fix Key Achievements to use proper\n# accomplishment statements\".\n\n$d = $word.ActiveDocument\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" section heading. There is\n# only one such heading in the document, so this anchors us reliably\n# even though some bullet text is duplicated elsewhere (e.g. under\n# \"Partner - Siege Analytics\").\n$paraCount = $d.Paragraphs.Count\n$sectionHeadingIndex = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionHeadingIndex = $i\n        break\n    }\n}\nif ($sectionHeadingIndex -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# Collect the bullet paragraph indices belonging to this section:\n# everything after the heading (and its optional \"Impact\" sub-heading)\n# up to (but excluding) the next heading-styled paragraph.\n$bulletStart = $sectionHeadingIndex + 1\nif ($d.Paragraphs.Item($bulletStart).Style.NameLocal -like \"Heading*\") {\n    $bulletStart = $bulletStart + 1\n}\n\n$bulletEnd = $bulletStart\nwhile (($bulletEnd -le $paraCount) -and ($d.Paragraphs.Item($bulletEnd).Style.NameLocal -notlike \"Heading*\")) {\n    $bulletEnd = $bulletEnd + 1\n}\n$bulletEnd = $bulletEnd - 1   # inclusive last bullet index\n\n$newBullets = @(\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    \"\u2022 `$4.7M savings enabled nonprofit access\",\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n    \"\u2022 178% accuracy improvement in racial classification algorithms\"\n)\n\n$bulletCount = $bulletEnd - $bulletStart + 1\nif ($bulletCount -lt $newBullets.Count) {\n    throw \"Fewer existing bullet paragraphs than replacement bullets\"\n}\n\n# Overwrite the text of the first N bullet paragraphs in place so their\n# formatting/paragraph identity is preserved.\nfor ($i = 0; $i -lt $newBullets.Count; $i++) {\n    $p = $d.Paragraphs.Item($bulletStart + $i)\n    $p.Range.Text = $newBullets[$i]\n}\n\n# Delete any leftover bullet paragraphs beyond the new bullet count,\n# working from the bottom up so indices of not-yet-deleted paragraphs\n# stay valid.\n$firstExtra = $bulletStart + $newBullets.Count\nfor ($i = $bulletEnd; $i -ge $firstExtra; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
